$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Formula = '=""&"63.506.98"'
$ws.Range("D2").Copy()
$ws.Range("D2").PasteSpecial(-4163)
$ws.Range("E2").Value = '  +1.37%  '
$ws.Range("D3").Formula = '=""&"2.650.80"'
$ws.Range("D3").Copy()
$ws.Range("D3").PasteSpecial(-4163)
$ws.Range("E3").Value = '  +2.74%  '
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").Formula = '=""&"591.24"'
$ws.Range("D5").Copy()
$ws.Range("D5").PasteSpecial(-4163)
$ws.Range("E5").Value = '  +1.72%  '
$ws.Range("D6").Formula = '=""&"143.96"'
$ws.Range("D6").Copy()
$ws.Range("D6").PasteSpecial(-4163)
$ws.Range("E6").Value = '  -0.72%  '
$ws.Range("E7").Value = '  +0.10%  '
$ws.Range("E8").Value = '  -0.92%  '
$ws.Range("D9").Formula = '=""&"2.649.36"'
$ws.Range("D9").Copy()
$ws.Range("D9").PasteSpecial(-4163)
$ws.Range("E9").Value = '  +2.74%  '
$ws.Range("E10").Value = '  +0.62%  '
$ws.Range("D11").Formula = '=""&"5.62"'
$ws.Range("D11").Copy()
$ws.Range("D11").PasteSpecial(-4163)
$ws.Range("E11").Value = '  +0.97%  '
$ws.Range("E12").Value = '  +0.84%  '
$ws.Range("E13").Value = '  +0.63%  '
$ws.Range("D14").Formula = '=""&"27.42"'
$ws.Range("D14").Copy()
$ws.Range("D14").PasteSpecial(-4163)
$ws.Range("E14").Value = '  +1.68%  '
$ws.Range("D15").Formula = '=""&"3.128.27"'
$ws.Range("D15").Copy()
$ws.Range("D15").PasteSpecial(-4163)
$ws.Range("E15").Value = '  +2.88%  '
$ws.Range("D16").Formula = '=""&"63.412.32"'
$ws.Range("D16").Copy()
$ws.Range("D16").PasteSpecial(-4163)
$ws.Range("E16").Value = '  +1.38%  '
$ws.Range("D17").Formula = '=""&"0.0000145"'
$ws.Range("D17").Copy()
$ws.Range("D17").PasteSpecial(-4163)
$ws.Range("E17").Value = '  +0.90%  '
$ws.Range("D18").Formula = '=""&"2.656.23"'
$ws.Range("D18").Copy()
$ws.Range("D18").PasteSpecial(-4163)
$ws.Range("E18").Value = '  +3.09%  '
$ws.Range("D19").Formula = '=""&"11.40"'
$ws.Range("D19").Copy()
$ws.Range("D19").PasteSpecial(-4163)
$ws.Range("E19").Value = '  +1.86%  '
$ws.Range("D20").Formula = '=""&"340.01"'
$ws.Range("D20").Copy()
$ws.Range("D20").PasteSpecial(-4163)
$ws.Range("E20").Value = '  +0.52%  '
$ws.Range("E21").Value = '  +0.69%  '
$ws.Range("D22").Formula = '=""&"6.73"'
$ws.Range("D22").Copy()
$ws.Range("D22").PasteSpecial(-4163)
$ws.Range("E22").Value = '  +1.26%  '
$ws.Range("E23").Value = '  +0.10%  '
$ws.Range("D24").Formula = '=""&"67.16"'
$ws.Range("D24").Copy()
$ws.Range("D24").PasteSpecial(-4163)
$ws.Range("E24").Value = '  +0.30%  '
$ws.Range("E25").Value = '  +5.20%  '
$ws.Range("E26").Value = '  +4.86%  '
$ws.Range("E27").Value = '  +0.55%  '
$ws.Range("D28").Formula = '=""&"542.40"'
$ws.Range("D28").Copy()
$ws.Range("D28").PasteSpecial(-4163)
$ws.Range("E28").Value = '  +17.63%  '
$ws.Range("E29").Value = '  +0.09%  '
$ws.Range("D30").Formula = '=""&"8.41"'
$ws.Range("D30").Copy()
$ws.Range("D30").PasteSpecial(-4163)
$ws.Range("E30").Value = '  +2.46%  '
$ws.Range("D31").Formula = '=""&"7.79"'
$ws.Range("D31").Copy()
$ws.Range("D31").PasteSpecial(-4163)
$ws.Range("E31").Value = '  -1.07%  '
$ws.Range("D32").Formula = '=""&"1.83"'
$ws.Range("D32").Copy()
$ws.Range("D32").PasteSpecial(-4163)
$ws.Range("E32").Value = '  +14.37%  '
$ws.Range("D33").Formula = '=""&"1.96"'
$ws.Range("D33").Copy()
$ws.Range("D33").PasteSpecial(-4163)
$ws.Range("E33").Value = '  +1.99%  '
$ws.Range("D34").Formula = '=""&"0.0₃0806"'
$ws.Range("D34").Copy()
$ws.Range("D34").PasteSpecial(-4163)
$ws.Range("E34").Value = '  -0.25%  '
$ws.Range("D35").Formula = '=""&"174.92"'
$ws.Range("D35").Copy()
$ws.Range("D35").PasteSpecial(-4163)
$ws.Range("E35").Value = '  -1.07%  '
$ws.Range("D36").Formula = '=""&"4.86"'
$ws.Range("D36").Copy()
$ws.Range("D36").PasteSpecial(-4163)
$ws.Range("E36").Value = '  +9.01%  '
$ws.Range("E37").Value = '  +0.06%  '
$ws.Range("E38").Value = '  +0.95%  '
$ws.Range("D39").Formula = '=""&"19.08"'
$ws.Range("D39").Copy()
$ws.Range("D39").PasteSpecial(-4163)
$ws.Range("E39").Value = '  +0.97%  '
$ws.Range("E40").Value = '  +7.11%  '
$ws.Range("D41").Formula = '=""&"170.81"'
$ws.Range("D41").Copy()
$ws.Range("D41").PasteSpecial(-4163)
$ws.Range("E41").Value = '  +8.64%  '
$ws.Range("D43").Formula = '=""&"40.27"'
$ws.Range("D43").Copy()
$ws.Range("D43").PasteSpecial(-4163)
$ws.Range("E43").Value = '  +2.17%  '
$ws.Range("B44").Value = 'InjectiveProtocol'
$ws.Range("C44").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D44").Formula = '=""&"22.45"'
$ws.Range("D44").Copy()
$ws.Range("D44").PasteSpecial(-4163)
$ws.Range("E44").Value = '  +6.25%  '
$ws.Range("B45").Value = 'Filecoin'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D45").Formula = '=""&"3.73"'
$ws.Range("D45").Copy()
$ws.Range("D45").PasteSpecial(-4163)
$ws.Range("E45").Value = '  +0.45%  '
$ws.Range("B46").Value = 'Hedera'
$ws.Range("C46").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D46").Formula = '=""&"0.0562"'
$ws.Range("D46").Copy()
$ws.Range("D46").PasteSpecial(-4163)
$ws.Range("E46").Value = '  +5.24%  '
$ws.Range("B47").Value = 'Mantle'
$ws.Range("C47").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D47").Formula = '=""&"0.632"'
$ws.Range("D47").Copy()
$ws.Range("D47").PasteSpecial(-4163)
$ws.Range("E47").Value = '  +0.79%  '
$ws.Range("E48").Value = '  -0.27%  '
$ws.Range("E49").Value = '  +2.27%  '
$ws.Range("D50").Formula = '=""&"18.79"'
$ws.Range("D50").Copy()
$ws.Range("D50").PasteSpecial(-4163)
$ws.Range("E50").Value = '  +3.93%  '
$ws.Range("E51").Value = '  +0.78%  '

$excel.CutCopyMode = 0
